$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.193.60'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.836.88'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9984'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.65'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6682'
$ws.Range('E6').Value = '  -2.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9994'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07402'
$ws.Range('E8').Value = '  -1.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2956'
$ws.Range('E9').Value = '  -2.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.88'
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07721'
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.837.89'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.022'
$ws.Range('E13').Value = '  -1.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6782'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '86.34'
$ws.Range('E15').Value = '  -3.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.195'
$ws.Range('E16').Value = '  -1.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008252'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '28.878.72'
$ws.Range('E18').Value = '  -1.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '228.90'
$ws.Range('E19').Value = '  -2.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.54'
$ws.Range('E20').Value = '  -0.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9987'
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.273'
$ws.Range('E22').Value = '  -3.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9993'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '160.21'
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.715'
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1407'
$ws.Range('E26').Value = '  -3.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.02'
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.506'
$ws.Range('E28').Value = '  -1.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.208'
$ws.Range('E29').Value = '  -0.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.090'
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.194'
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05359'
$ws.Range('E32').Value = '  +3.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.871'
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7507'
$ws.Range('E34').Value = '  -2.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.139'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.677'
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.330.18'
$ws.Range('E37').Value = '  +2.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01802'
$ws.Range('E38').Value = '  -2.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.735'
$ws.Range('E39').Value = '  +1.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9211'
$ws.Range('E40').Value = '  -2.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.969'
$ws.Range('E41').Value = '  +4.14%  '
$ws.Range('B42').Value = 'XinFinNetwork'
$ws.Range('C42').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.08346'
$ws.Range('E42').Value = '  +20.79%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9992'
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '103.23'
$ws.Range('E44').Value = '  -2.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000124'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5170'
$ws.Range('E46').Value = '  -0.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.955.42'
$ws.Range('E47').Value = '  -1.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '64.00'
$ws.Range('E48').Value = '  +1.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.762'
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.227'
$ws.Range('E50').Value = '  -5.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05944'
$ws.Range('E51').Value = '  +0.19%  '
